# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 21:22"

# --- Update totals for Estados Unidos (row 4) ---
$ws.Range("B4").Value = 868482
$ws.Range("C4").Value = 19765
$ws.Range("D4").Value = 84825
$ws.Range("E4").Value = 734563
$ws.Range("G4").Value = 1435
$ws.Range("H4").Value = 49094

# --- Update totals for Peru (row 21) ---
$ws.Range("D21").Value = 7422
$ws.Range("E21").Value = 12920

# --- Update totals for Costa Rica (row 96) ---
$ws.Range("B96").Value = 687
$ws.Range("C96").Value = 6
$ws.Range("D96").Value = 196
$ws.Range("E96").Value = 485
$ws.Range("F96").Value = 8

# --- Update totals for Georgia (row 110) ---
$ws.Range("B110").Value = 425
$ws.Range("C110").Value = 9
$ws.Range("E110").Value = 309

# --- Reorder Macao / Siria / Mozambique rows (166-168) and refresh their
#     totals, matching the new data pull where Mozambique overtakes Macao ---
$ws.Range("A166").Value = "Mozambique"
$ws.Range("B166").Value = 46
$ws.Range("C166").Value = 5
$ws.Range("D166").Value = 9
$ws.Range("E166").Value = 37
$ws.Range("F166").Value = 0

$ws.Range("A167").Value = "Macao"
$ws.Range("B167").Value = 45
$ws.Range("D167").Value = 27
$ws.Range("E167").Value = 18
$ws.Range("F167").Value = 1
$ws.Range("H167").Value = 0

$ws.Range("A168").Value = "Siria"
$ws.Range("B168").Value = 42
$ws.Range("D168").Value = 6
$ws.Range("H168").Value = 3
